# Config Workflow.xlsx - "gamble fe10 to fe15 again"
# Resets the FE10..FE15 hyperparameter gamble values on sheet "FE",
# resets the "FE10" status on "Sheet1" back to "corriendo", and
# updates sheet selections / the active tab to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: FE10 status cell goes back from "ok" to "corriendo", and the
# selection on that sheet moves from E16:E18 to just B16.
# ---------------------------------------------------------------------
$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Range("E16").Value = "corriendo"

# ---------------------------------------------------------------------
# FE: re-roll (gamble) the num.trees / max.depth / min.node.size
# hyperparameters for rows 12-17 (FE10..FE15).
# ---------------------------------------------------------------------
$wsFE = $wb.Worksheets.Item("FE")

$wsFE.Range("B12").Value = 200
$wsFE.Range("C12").Value = 10
$wsFE.Range("D12").Value = 300

$wsFE.Range("B13").Value = 200
$wsFE.Range("C13").Value = 10
$wsFE.Range("D13").Value = 400

$wsFE.Range("B14").Value = 200
$wsFE.Range("C14").Value = 10
$wsFE.Range("D14").Value = 600

$wsFE.Range("B15").Value = 200
$wsFE.Range("C15").Value = 10
$wsFE.Range("D15").Value = 800

$wsFE.Range("B16").Value = 200
$wsFE.Range("C16").Value = 10
$wsFE.Range("D16").Value = 1000

$wsFE.Range("B17").Value = 200
$wsFE.Range("C17").Value = 10
$wsFE.Range("D17").Value = 1200

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping: the FE sheet becomes the active
# (selected) tab, with B12:B17 selected; Sheet1 ends up with B16
# selected; TS loses its tabSelected flag since FE is now active.
# ---------------------------------------------------------------------
$wsSheet1.Activate()
$wsSheet1.Range("B16").Select()

$wsFE.Activate()
$wsFE.Range("B12:B17").Select()
